$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.100.72"
$ws.Range("D3").Value = "1.656.67"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5305"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2616"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06338"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07764"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.494"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").Value = "1.650.89"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5477"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "0.0₅8153"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "26.127.65"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.550"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.011"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "140.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1244"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.286"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05948"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.276"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.513"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.240"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.547"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9470"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.29%  "
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5640"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.27%  "
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.847"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8482"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.002"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "1.014.22"
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").Value = "1.801.08"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").Value = "0.0₈105"
$ws.Range("E46").Value = "  -6.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.004"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05154"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.470"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.753"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.00%  "
